$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the analysis timestamp string (shared string used by A2)
$ws.Range("A2").Value = "2025-05-30 10:33:41"

# Update the metrics values for row 2 per the new extraction
$values = @{
    "B2"  = 2550
    "C2"  = 1436
    "D2"  = 56.31372549019608
    "E2"  = 347
    "F2"  = 13.6078431372549
    "G2"  = 1718
    "H2"  = 67.37254901960785
    "I2"  = 811
    "J2"  = 31.80392156862745
    "K2"  = 352934.79
    "L2"  = 21
    "M2"  = 0.823529411764706
    "N2"  = 1110.97
    "O2"  = 84
    "P2"  = 3.294117647058824
    "Q2"  = 7248.179999999999
    "R2"  = 543
    "S2"  = 21.29411764705883
    "T2"  = 696
    "U2"  = 27.29411764705883
    "V2"  = 342393.66
    "W2"  = 460
    "X2"  = 18.03921568627451
    "Y2"  = 31
    "Z2"  = 1.215686274509804
    "AA2" = 3292.95
    "AB2" = 736
    "AC2" = 28.86274509803922
    "AD2" = 2550
    "AE2" = 2323
    "AF2" = 91.09803921568627
    "AG2" = 227
    "AH2" = 8.901960784313729
    "AI2" = 129
    "AJ2" = 179
    "AK2" = 364
    "AL2" = 19.19642857142857
    "AM2" = 26.63690476190476
    "AN2" = 54.16666666666666
    "AO2" = 487005.88
    "AP2" = 92306.54000000001
    "AQ2" = 30577.61
    "AR2" = 79.85142501837585
    "AS2" = 15.13494818074006
    "AT2" = 5.013626800884087
    "AU2" = 77.1710575509085
    "AV2" = 116.6083640985293
    "AW2" = 156.8895800933126
}

foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}
